$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191, shifting existing rows 191:219 down to 192:220.
$ws.Rows("191:191").Insert()

# Populate the newly inserted row 191 with the new weekly price record.
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 45124
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100114002
$ws.Range("G191").Value = "Camote"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 100
$ws.Range("K191").Value = 26000
$ws.Range("L191").Value = 26000
$ws.Range("M191").Value = 26000
$ws.Range("N191").Value = "$/caja 18 kilos"
$ws.Range("O191").Value = "Perú"
$ws.Range("P191").Value = 1444
$ws.Range("Q191").Value = 18
$ws.Range("R191").Value = "Hortaliza"
